# The document contains a merge-field whose field code is:
#     { m:'doc.html'.fromHTMLURI() }
# built out of <w:fldChar>/<w:instrText> runs. The parser that consumes
# this template was switched to TokenIteratorFieldRewriterSplit, which
# expects the token to be expressed as plain literal text runs instead
# of a real Word field, i.e.:
#     {m:'doc.html'.fromHTMLURI()}
# split across <w:t> runs, keeping the (hidden) "_GoBack" bookmark in
# place in the middle of the run sequence.

$d = $word.ActiveDocument

# Locate the field and the paragraph that hosts it (robust to paragraph
# numbering: look up the paragraph whose range contains the field code).
$field = $d.Fields.Item(1)
$codeStart = $field.Code.Start

$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($codeStart -ge $candidate.Range.Start -and $codeStart -le $candidate.Range.End) {
        $targetParagraph = $candidate
        break
    }
}

$q = "'"

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>' + $q + '</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>' + $q + '.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

$targetParagraph.Range.InsertXML($newParagraphXml)
